# Update the FilesTab's TabQuery (B5) with the new SQL query text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newFilesQuery = @'
SELECT Distinct 
    seq.file_name AS "File Name",
    seq.data_category AS "Data Category",
    COALESCE(seq.file_description, '') AS "File Description",
    seq.file_type AS "File Type",
    CASE     
        WHEN seq.file_size >= 1024 * 1024 * 1024 THEN 
            CASE 
                WHEN ROUND(seq.file_size / (1024.0 * 1024.0 * 1024.0), 2) = 
                     CAST(ROUND(seq.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
                THEN CAST(CAST(ROUND(seq.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
                ELSE ROUND(seq.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
            END
        WHEN seq.file_size >= 1024 * 1024 THEN 
            CASE 
                WHEN ROUND(seq.file_size / (1024.0 * 1024.0), 2) = 
                     CAST(ROUND(seq.file_size / (1024.0 * 1024.0), 0) AS INT) 
                THEN CAST(CAST(ROUND(seq.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
                ELSE ROUND(seq.file_size / (1024.0 * 1024.0), 2) || ' MB'
            END
        WHEN seq.file_size >= 1024 THEN 
            CASE 
                WHEN ROUND(seq.file_size / 1024.0, 2) = 
                     CAST(ROUND(seq.file_size / 1024.0, 0) AS INT) 
                THEN CAST(CAST(ROUND(seq.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
                ELSE ROUND(seq.file_size / 1024.0, 2) || ' KB'
            END
        ELSE 
            CASE 
                WHEN ROUND(seq.file_size, 2) = 
                     CAST(ROUND(seq.file_size, 0) AS INT) 
                THEN CAST(CAST(ROUND(seq.file_size, 0) AS INT) AS TEXT) || ' Bytes'
                ELSE ROUND(seq.file_size, 2) || ' Bytes'
            END
    END AS "File Size",
    seq.file_access AS "File Access",
    std.dbgap_accession AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID"    
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"  
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
                          AND seq.library_source_molecule = 'Genomic'
                          AND seq.library_source_material = 'Single-cells'
WHERE 
    std.dbgap_accession = 'phs003519' 
    AND std.study_status = 'Completed'
    AND seq."sample.id" IS NOT NULL
ORDER BY  seq.file_name asc
LIMIT 100;
'@

$ws.Cells.Item(5, 2).Value = $newFilesQuery

# Re-assert the row height - writing a long wrapped string recalculates
# the autofit height; put it back to match the row's prior (capped) height.
$ws.Rows.Item(5).RowHeight = 409.6

# Reflect the latest on-screen selection/scroll position.
[void]$ws.Range("C10").Select()
